$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.214.54'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.904.76'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'307.82"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D7").Value = "'0.5264"
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("D8").Value = "'0.3824"
$ws.Range("E8").Value = '  +1.54%  '
$ws.Range("D9").Value = "'0.07302"
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("D10").Value = "'21.54"
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("E12").Value = '  -4.47%  '
$ws.Range("D13").Value = "'95.84"
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("D14").Value = "'5.371"
$ws.Range("E14").Value = '  +1.55%  '
$ws.Range("D15").Value = '1.801.05'
$ws.Range("E15").Value = '  -5.37%  '
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").Value = "'0.000008678"
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = "'14.73"
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").Value = '27.249.67'
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = "'5.122"
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").Value = "'10.84"
$ws.Range("E22").Value = '  +2.24%  '
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("D24").Value = "'2.339"
$ws.Range("E24").Value = '  +2.46%  '
$ws.Range("D25").Value = "'150.12"
$ws.Range("E25").Value = '  +1.92%  '
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("D27").Value = "'1.743"
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").Value = "'116.74"
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("D29").Value = "'4.846"
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = "'4.870"
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").Value = "'0.09230"
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").Value = "'0.8158"
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = "'1.229"
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("D35").Value = "'2.988"
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("D36").Value = "'3.360"
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").Value = "'2.705"
$ws.Range("E37").Value = '  +3.43%  '
$ws.Range("D38").Value = "'0.5729"
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = "'0.01996"
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("D41").Value = "'8.996"
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("D42").Value = "'6.627"
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = "'117.01"
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("D44").Value = "'0.1522"
$ws.Range("E44").Value = '  +0.64%  '
$ws.Range("D45").Value = "'0.4925"
$ws.Range("E45").Value = '  +1.30%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'10.21"
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").Value = "'38.56"
$ws.Range("E49").Value = '  +2.92%  '
$ws.Range("D50").Value = "'64.23"
$ws.Range("E50").Value = '  +0.34%  '

# Reset style on forced-text numeric cells to avoid quotePrefix/number-format artifacts
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
